$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- L250i price update (row 10)
$ws.Range("G10").Value = 1180
$ws.Range("H10").Value = 1270

# --- Right-hand price table (F:I, rows 18-41): a new model "T130" is inserted
# in alphabetical order before "T140", shifting all following rows down by one.
$ws.Range("F18").Value = "T130"
$ws.Range("G18").Value = 1250
$ws.Range("H18").Value = 1350
$ws.Range("I18").ClearContents()

# --- Footer note text (F4): "Murad = 01740883469" -> "Last Update: 10.11.2020"
$ws.Range("F4").Value = "Last Update: 10.11.2020"

$ws.Range("F19").Value = "T140"
$ws.Range("G19").Value = 1370
$ws.Range("H19").Value = 1490
$ws.Range("I19").ClearContents()

$ws.Range("F20").Value = "T180"
$ws.Range("G20").Value = 1220
$ws.Range("H20").Value = 1320
$ws.Range("I20").ClearContents()

$ws.Range("F21").Value = "V102_SKD"
$ws.Range("G21").Value = 3610
$ws.Range("H21").Value = 3890
$ws.Range("I21").ClearContents()

$ws.Range("F22").Value = "V105_SKD"
$ws.Range("G22").Value = 3890
$ws.Range("H22").Value = 4190
$ws.Range("I22").ClearContents()

$ws.Range("F23").Value = "V128_SKD"
$ws.Range("G23").Value = 4280
$ws.Range("H23").Value = 4590
$ws.Range("I23").ClearContents()

$ws.Range("F24").Value = "V141_SKD"
$ws.Range("G24").Value = 4180
$ws.Range("H24").Value = 4490
$ws.Range("I24").ClearContents()

$ws.Range("F25").Value = "V155"
$ws.Range("G25").Value = 5390
$ws.Range("H25").Value = 5790
$ws.Range("I25").ClearContents()

$ws.Range("F26").Value = "V44"
$ws.Range("G26").Value = 3560
$ws.Range("H26").Value = 3840
$ws.Range("I26").ClearContents()

$ws.Range("F27").Value = "V48_SKD"
$ws.Range("G27").Value = 3340
$ws.Range("H27").Value = 3590
$ws.Range("I27").ClearContents()

$ws.Range("F28").Value = "V75_SKD"
$ws.Range("G28").Value = 4500
$ws.Range("H28").Value = 4790
$ws.Range("I28").ClearContents()

$ws.Range("F29").Value = "V94_SKD"
$ws.Range("G29").Value = 3620
$ws.Range("H29").Value = 3890
$ws.Range("I29").ClearContents()

$ws.Range("F30").Value = "V97_SKD"
$ws.Range("G30").Value = 4080
$ws.Range("H30").Value = 4390
$ws.Range("I30").ClearContents()

$ws.Range("F31").Value = "V98_SKD"
$ws.Range("G31").Value = 4220
$ws.Range("H31").Value = 4540
$ws.Range("I31").ClearContents()

$ws.Range("F32").Value = "V99 +_SKD"
$ws.Range("G32").Value = 3640
$ws.Range("H32").Value = 3890
$ws.Range("I32").ClearContents()

$ws.Range("F33").Value = "V99_SKD"
$ws.Range("G33").Value = 3710
$ws.Range("H33").Value = 3990
$ws.Range("I33").ClearContents()

$ws.Range("F34").Value = "Z12_SKD"
$ws.Range("G34").Value = 7350
$ws.Range("H34").Value = 7990
$ws.Range("I34").Value = 700

$ws.Range("F35").Value = "Z15_SKD"
$ws.Range("G35").Value = 7890
$ws.Range("H35").Value = 8490
$ws.Range("I35").Value = 300

$ws.Range("F36").Value = "Z16_SKD"
$ws.Range("G36").Value = 7790
$ws.Range("H36").Value = 8290
$ws.Range("I36").ClearContents()

$ws.Range("F37").Value = "Z20_SKD"
$ws.Range("G37").Value = 8310
$ws.Range("H37").Value = 8990
$ws.Range("I37").ClearContents()

$ws.Range("F38").Value = "Z25_SKD"
$ws.Range("G38").Value = 7800
$ws.Range("H38").Value = 8390
$ws.Range("I38").ClearContents()

$ws.Range("F39").Value = "Z28_SKD"
$ws.Range("G39").Value = 8450
$ws.Range("H39").Value = 8990
$ws.Range("I39").ClearContents()

$ws.Range("F40").Value = "Z30_SKD"
$ws.Range("G40").Value = 9300
$ws.Range("H40").Value = 9790
$ws.Range("I40").ClearContents()

$ws.Range("F41").Value = "Z50_SKD"
$ws.Range("G41").Value = 10340
$ws.Range("H41").Value = 10990
$ws.Range("I41").ClearContents()

# --- The highlighted row's fill/format moves down from row 37 to row 38
# (it always tracks the physical row, independent of the data it holds).
$normalRange = $ws.Range("F37:I37")
$normalRange.Interior.Pattern = -4142
$ws.Range("F37").Font.Bold = $true
$ws.Range("G37:H37").Font.Bold = $false
$ws.Range("I37").Font.Bold = $false
$ws.Range("G37:H37").NumberFormat = "#,##0"

$highlightRange = $ws.Range("F38:I38")
$highlightRange.Interior.Color = 65535
$ws.Range("F38").Font.Bold = $true
$ws.Range("G38:H38").Font.Bold = $false
$ws.Range("I38").Font.Bold = $false
$ws.Range("G38:H38").NumberFormat = "#,##0"

# --- Active cell selection moved from L9 to L10
[void]$ws.Range("L10").Select()
